# Applies the "update scripts with new tpm" data refresh to the
# NatMI LR-pairs sheet (Agrp-Sdc3), recomputing expression/specificity
# metrics for the updated TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5679715
$ws.Range("H2").Value = 1.135943
$ws.Range("I2").Value = 0.1948020679094191
$ws.Range("J2").Value = 0.1388865047139418
$ws.Range("M2").Value = 35.18694
$ws.Range("N2").Value = 70.37388
$ws.Range("O2").Value = 0.1785572969025014
$ws.Range("P2").Value = 0.1282643028201015
$ws.Range("Q2").Value = 19.98517909221
$ws.Range("R2").Value = 79.94071636884
$ws.Range("S2").Value = 0.03478333067692338
$ws.Range("T2").Value = 0.01781418069825448

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5679715
$ws.Range("H3").Value = 1.135943
$ws.Range("I3").Value = 0.1948020679094191
$ws.Range("J3").Value = 0.1388865047139418
$ws.Range("O3").Value = 0.02719713854783473
$ws.Range("P3").Value = 0.02930506404712654
$ws.Range("Q3").Value = 3.0440631332525
$ws.Range("R3").Value = 18.264378799515
$ws.Range("S3").Value = 0.005298058830337181
$ws.Range("T3").Value = 0.004070077915923607

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5679715
$ws.Range("H4").Value = 1.135943
$ws.Range("I4").Value = 0.1948020679094191
$ws.Range("J4").Value = 0.1388865047139418
$ws.Range("M4").Value = 70.59161999999999
$ws.Range("N4").Value = 211.77486
$ws.Range("O4").Value = 0.3582195226742806
$ws.Range("P4").Value = 0.3859834752997077
$ws.Range("Q4").Value = 40.09402829882999
$ws.Range("R4").Value = 240.56416979298
$ws.Range("S4").Value = 0.06978190378247492
$ws.Range("T4").Value = 0.05360789576171649

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5679715
$ws.Range("H5").Value = 1.135943
$ws.Range("I5").Value = 0.1948020679094191
$ws.Range("J5").Value = 0.1388865047139418
$ws.Range("M5").Value = 7.337415
$ws.Range("N5").Value = 14.67483
$ws.Range("O5").Value = 0.03723395636710288
$ws.Range("P5").Value = 0.0267465263952124
$ws.Range("Q5").Value = 4.1674426036725
$ws.Range("R5").Value = 16.66977041469
$ws.Range("S5").Value = 0.007253251696760724
$ws.Range("T5").Value = 0.003714731564270236

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5679715
$ws.Range("H6").Value = 1.135943
$ws.Range("I6").Value = 0.1948020679094191
$ws.Range("J6").Value = 0.1388865047139418
$ws.Range("M6").Value = 15.36873766666667
$ws.Range("N6").Value = 46.106213
$ws.Range("O6").Value = 0.07798917025929666
$ws.Range("P6").Value = 0.0840337532351506
$ws.Range("Q6").Value = 8.729004985643165
$ws.Range("R6").Value = 52.37402991385899
$ws.Range("S6").Value = 0.01519245164105076
$ws.Range("T6").Value = 0.01167115426482396

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5679715
$ws.Range("H7").Value = 1.135943
$ws.Range("I7").Value = 0.1948020679094191
$ws.Range("J7").Value = 0.1388865047139418
$ws.Range("M7").Value = 63.21821133333334
$ws.Range("N7").Value = 189.654634
$ws.Range("O7").Value = 0.3208029152489838
$ws.Range("P7").Value = 0.3456668782027013
$ws.Range("Q7").Value = 35.90614231831033
$ws.Range("R7").Value = 215.436853909862
$ws.Range("S7").Value = 0.06249307128187217
$ws.Range("T7").Value = 0.04800846450895302

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.235202
$ws.Range("H8").Value = 6.705606
$ws.Range("I8").Value = 0.7666264448044829
$ws.Range("J8").Value = 0.8198634784745682
$ws.Range("M8").Value = 35.18694
$ws.Range("N8").Value = 70.37388
$ws.Range("O8").Value = 0.1785572969025014
$ws.Range("P8").Value = 0.1282643028201015
$ws.Range("Q8").Value = 78.64991866187999
$ws.Range("R8").Value = 471.89951197128
$ws.Range("S8").Value = 0.1368867457182631
$ws.Range("T8").Value = 0.1051592174742038

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.235202
$ws.Range("H9").Value = 6.705606
$ws.Range("I9").Value = 0.7666264448044829
$ws.Range("J9").Value = 0.8198634784745682
$ws.Range("O9").Value = 0.02719713854783473
$ws.Range("P9").Value = 0.02930506404712654
$ws.Range("Q9").Value = 11.97964335107
$ws.Range("R9").Value = 107.81679015963
$ws.Range("S9").Value = 0.0208500456337815
$ws.Range("T9").Value = 0.02402615174659718

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.235202
$ws.Range("H10").Value = 6.705606
$ws.Range("I10").Value = 0.7666264448044829
$ws.Range("J10").Value = 0.8198634784745682
$ws.Range("M10").Value = 70.59161999999999
$ws.Range("N10").Value = 211.77486
$ws.Range("O10").Value = 0.3582195226742806
$ws.Range("P10").Value = 0.3859834752997077
$ws.Range("Q10").Value = 157.78653020724
$ws.Range("R10").Value = 1420.07877186516
$ws.Range("S10").Value = 0.2746205591273426
$ws.Range("T10").Value = 0.316453754692921

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.235202
$ws.Range("H11").Value = 6.705606
$ws.Range("I11").Value = 0.7666264448044829
$ws.Range("J11").Value = 0.8198634784745682
$ws.Range("M11").Value = 7.337415
$ws.Range("N11").Value = 14.67483
$ws.Range("O11").Value = 0.03723395636710288
$ws.Range("P11").Value = 0.0267465263952124
$ws.Range("Q11").Value = 16.40060468283
$ws.Range("R11").Value = 98.40362809697999
$ws.Range("S11").Value = 0.02854453559571732
$ws.Range("T11").Value = 0.02192850016749069

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.235202
$ws.Range("H12").Value = 6.705606
$ws.Range("I12").Value = 0.7666264448044829
$ws.Range("J12").Value = 0.8198634784745682
$ws.Range("M12").Value = 15.36873766666667
$ws.Range("N12").Value = 46.106213
$ws.Range("O12").Value = 0.07798917025929666
$ws.Range("P12").Value = 0.0840337532351506
$ws.Range("Q12").Value = 34.35223317000866
$ws.Range("R12").Value = 309.170098530078
$ws.Range("S12").Value = 0.05978856032913612
$ws.Range("T12").Value = 0.06889620523664407

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.235202
$ws.Range("H13").Value = 6.705606
$ws.Range("I13").Value = 0.7666264448044829
$ws.Range("J13").Value = 0.8198634784745682
$ws.Range("M13").Value = 63.21821133333334
$ws.Range("N13").Value = 189.654634
$ws.Range("O13").Value = 0.3208029152489838
$ws.Range("P13").Value = 0.3456668782027013
$ws.Range("Q13").Value = 141.3054724086893
$ws.Range("R13").Value = 1271.749251678204
$ws.Range("S13").Value = 0.2459359984002423
$ws.Range("T13").Value = 0.2833996491567116

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.1124603333333333
$ws.Range("H14").Value = 0.337381
$ws.Range("I14").Value = 0.03857148728609783
$ws.Range("J14").Value = 0.04125001681149001
$ws.Range("M14").Value = 35.18694
$ws.Range("N14").Value = 70.37388
$ws.Range("O14").Value = 0.1785572969025014
$ws.Range("P14").Value = 0.1282643028201015
$ws.Range("Q14").Value = 3.95713500138
$ws.Range("R14").Value = 23.74281000828
$ws.Range("S14").Value = 0.006887220507314826
$ws.Range("T14").Value = 0.005290904647643232

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.1124603333333333
$ws.Range("H15").Value = 0.337381
$ws.Range("I15").Value = 0.03857148728609783
$ws.Range("J15").Value = 0.04125001681149001
$ws.Range("O15").Value = 0.02719713854783473
$ws.Range("P15").Value = 0.02930506404712654
$ws.Range("Q15").Value = 0.6027350926116667
$ws.Range("R15").Value = 5.424615833504999
$ws.Range("S15").Value = 0.001049034083716048
$ws.Range("T15").Value = 0.001208834384605762

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.1124603333333333
$ws.Range("H16").Value = 0.337381
$ws.Range("I16").Value = 0.03857148728609783
$ws.Range("J16").Value = 0.04125001681149001
$ws.Range("M16").Value = 70.59161999999999
$ws.Range("N16").Value = 211.77486
$ws.Range("O16").Value = 0.3582195226742806
$ws.Range("P16").Value = 0.3859834752997077
$ws.Range("Q16").Value = 7.938757115739999
$ws.Range("R16").Value = 71.44881404166
$ws.Range("S16").Value = 0.01381705976446305
$ws.Range("T16").Value = 0.01592182484507029

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.1124603333333333
$ws.Range("H17").Value = 0.337381
$ws.Range("I17").Value = 0.03857148728609783
$ws.Range("J17").Value = 0.04125001681149001
$ws.Range("M17").Value = 7.337415
$ws.Range("N17").Value = 14.67483
$ws.Range("O17").Value = 0.03723395636710288
$ws.Range("P17").Value = 0.0267465263952124
$ws.Range("Q17").Value = 0.825168136705
$ws.Range("R17").Value = 4.95100882023
$ws.Range("S17").Value = 0.00143616907462483
$ws.Range("T17").Value = 0.001103294663451473

# Row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.1124603333333333
$ws.Range("H18").Value = 0.337381
$ws.Range("I18").Value = 0.03857148728609783
$ws.Range("J18").Value = 0.04125001681149001
$ws.Range("M18").Value = 15.36873766666667
$ws.Range("N18").Value = 46.106213
$ws.Range("O18").Value = 0.07798917025929666
$ws.Range("P18").Value = 0.0840337532351506
$ws.Range("Q18").Value = 1.728373360905889
$ws.Range("R18").Value = 15.555360248153
$ws.Range("S18").Value = 0.00300815828910978
$ws.Range("T18").Value = 0.003466393733682566

# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.1124603333333333
$ws.Range("H19").Value = 0.337381
$ws.Range("I19").Value = 0.03857148728609783
$ws.Range("J19").Value = 0.04125001681149001
$ws.Range("M19").Value = 63.21821133333334
$ws.Range("N19").Value = 189.654634
$ws.Range("O19").Value = 0.3208029152489838
$ws.Range("P19").Value = 0.3456668782027013
$ws.Range("Q19").Value = 7.109541119283778
$ws.Range("R19").Value = 63.985870073554
$ws.Range("S19").Value = 0.0123738455668693
$ws.Range("T19").Value = 0.0142587645370367

